# Apply the "feat: add 2022-Q3 data" change:
#  1. Update the "总计" (totals) sheet: insert a new first data row for
#     2022-Q3 and push the existing 2021-Q4 / 2020-Q4 rows down by one.
#  2. Insert a new "2022-Q3" worksheet between "总计" and "2021-Q4",
#     built from a copy of the "2021-Q4" sheet (same layout/styles),
#     populated with the 2022-Q3 fund holdings.

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Update the "总计" summary sheet.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Push the existing rows down first (bottom-up so we don't clobber data
# before it has been copied/re-used).
$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial($xlPasteFormats)
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2020-Q4"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.08

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 7
$total.Range("D3").Value = 1.78

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 8
$total.Range("D2").Value = 0.09

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q3" sheet (cloned from "2021-Q4" so it keeps
#    identical headers/styling), positioned right before "2021-Q4".
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($q4)
$new = $wb.Worksheets.Item("2021-Q4 (2)")
$new.Name = "2022-Q3"

# The template only has 7 data rows (r2:r8); we need an 8th (r9). Clone
# the row-8 styling down to row 9 before filling in values.
$new.Range("A8").Copy()
$new.Range("A9").PasteSpecial($xlPasteFormats)

# Fund holdings data for 2022-Q3 (A=index, B..G text fields, H numeric rank).
$data = @(
    @(0, "011097", "达诚宜创精选混合A", "0.70", "73.65", "3.15", "0.0220", 6),
    @(1, "010301", "达诚成长先锋混合A", "0.52", "74.13", "3.15", "0.0164", 5),
    @(2, "010808", "达诚策略先锋混合A", "0.33", "75.14", "3.19", "0.0105", 6),
    @(3, "010809", "达诚策略先锋混合C", "0.32", "75.14", "3.19", "0.0102", 6),
    @(4, "010302", "达诚成长先锋混合C", "0.31", "74.13", "3.15", "0.0098", 5),
    @(5, "011031", "达诚价值先锋灵活配置混合C", "0.23", "75.36", "3.27", "0.0075", 7),
    @(6, "011030", "达诚价值先锋灵活配置混合A", "0.21", "75.36", "3.27", "0.0069", 7),
    @(7, "011098", "达诚宜创精选混合C", "0.18", "73.65", "3.15", "0.0057", 6)
)

# Pre-format the text columns as Text so codes like "011097" and decimal
# strings like "0.70" survive as literal text (matching the source data,
# which stores them as inline strings rather than numbers) instead of
# being auto-coerced to numbers (and losing significant trailing zeros).
$new.Range("B2:G9").NumberFormat = "@"

$textCols = @("B", "C", "D", "E", "F", "G")
$r = 2
foreach ($row in $data) {
    $new.Range("A$r").Value = $row[0]
    for ($i = 1; $i -le 6; $i++) {
        $col = $textCols[$i - 1]
        $new.Range("$col$r").Value = $row[$i]
    }
    $new.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# The values are now committed as text; drop the temporary Text number
# format again so the cells end up with the same (default) styling as
# the rest of the sheet.
$new.Range("B2:G9").ClearFormats()

# Restore the originally-selected sheet/cell state (2020-Q4 was the
# active tab before the edit).
$wb.Worksheets.Item("2020-Q4").Select()
